#
# Applies the "第13组/13项目计划表.xlsx" update:
#  - adds a new weekly block (date header + table + summary) for
#    "2018.11.05 第十周周三" as rows 135-143 on Sheet1
#  - fills in the previously-empty completion % for the prior block
#    (rows 127-131, column C)
#  - keeps all pre-existing formatting untouched by cloning it (via
#    Copy + PasteSpecial formats) from the nearest matching template
#    rows already in the sheet, instead of inventing brand-new styles
#
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------
# 1. Completion-rate values that were left blank in the previous
#    week's block (rows 127-131, column C)
# ---------------------------------------------------------------
$ws.Range("C127").Value = 1
$ws.Range("C128").Value = 0.9
$ws.Range("C129").Value = 1
$ws.Range("C130").Value = 1
$ws.Range("C131").Value = 1

# ---------------------------------------------------------------
# 2. Merge the header / summary rows exactly like every other week
#    block in the sheet (A:D on the header row and the two "总结："
#    rows). Doing this BEFORE cloning formatting keeps the border
#    table from being split into per-cell fragments the way it
#    would if Merge() ran after the cells already had the thin-box
#    border applied.
# ---------------------------------------------------------------
$ws.Range("A135:D135").Merge()
$ws.Range("A142:D143").Merge()

# ---------------------------------------------------------------
# 3. Clone the row formatting for the new block from the closest
#    existing analogues so the new rows match the sheet's visual
#    style exactly (same borders / fonts / number formats /
#    alignment as every other week block).
# ---------------------------------------------------------------
# Row 135 = new date-header row -> looks like row 31's header
$ws.Range("A31:D31").Copy()
$ws.Range("A135:D135").PasteSpecial(-4122)

# Row 136 (column titles) + rows 137-141 (data rows) -> look like
# rows 126-131 (the immediately preceding table)
$ws.Range("A126:D131").Copy()
$ws.Range("A136:D141").PasteSpecial(-4122)

# Rows 142-143 ("总结：" summary rows) -> look like rows 132-133
$ws.Range("A132:D133").Copy()
$ws.Range("A142:D143").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# 3. Cell values / text for the new block
# ---------------------------------------------------------------
$ws.Range("A135").Value = "日期：2018.11.05 第十周周三"

$ws.Range("A136").Value = "组员"
$ws.Range("B136").Value = "计划内容"
$ws.Range("C136").Value = "完成情况"
$ws.Range("D136").Value = "备注"

$ws.Range("A137").Value = "王伟锋"
$ws.Range("B137").Value = "完成群组的所有后台编写"

$ws.Range("A138").Value = "陈升云"
$ws.Range("B138").Value = "完成个人资料的查看，修改等"

$ws.Range("A139").Value = "林玮成"
$ws.Range("B139").Value = "辅助app开发"

$ws.Range("A140").Value = "吴帅辰"
$ws.Range("B140").Value = "基本功能已完成，现优化管理员系统"

$ws.Range("A141").Value = "李海洋"
$ws.Range("B141").Value = "获取群成员的定位信息"

$ws.Range("A142").Value = "总结："

# ---------------------------------------------------------------
# 4. Match the saved selection / scroll position of the edited
#    workbook (cursor left on D140 after typing the last entry).
# ---------------------------------------------------------------
$ws.Range("D140").Select()
